$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Hide the previously-reviewed example rows (2-6) -----------------------
for ($r = 2; $r -le 6; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# --- Row 7: new plain-text sentence example + its source URL ---------------
$ws.Range("A7").Value = "Europska komisija je bi sredinom sljedećeg mjeseca trebala dati pozitivnu ocjenu da je Hrvatska ispunila sve kriterije za ulazak u šengenski prostor"
$ws.Range("B7").Value = "https://www.24sata.hr/news/bitno-je-dobiti-zeleno-svjetlo-ulazak-ne-mozemo-predvidjeto-650975"

# --- Row 8: sentence with one bold run + its source URL ---------------------
$a8Prefix = "U tom pravcu, nama je interes da se "
$a8Bold   = "hladne glave"
$a8Suffix = " pristupi svemu."
$ws.Range("A8").Value = $a8Prefix + $a8Bold + $a8Suffix

$pos = 1
$len = $a8Prefix.Length
$pos += $len

$len = $a8Bold.Length
$run = $ws.Range("A8").Characters($pos, $len)
$run.Font.Bold = $true
$run.Font.Name = "Calibri"
$run.Font.Size = 12
$pos += $len

$len = $a8Suffix.Length
$run = $ws.Range("A8").Characters($pos, $len)
$run.Font.Name = "Calibri"
$run.Font.Size = 12

$ws.Range("B8").Value = "http://novilist.hr/Vijesti/Hrvatska/Kristianova-majka-ne-zeli-da-se-ikome-vise-ponovi-ono-sto-je-ubilo-njenog-sina-Trazimo-pravdu"

# --- Row 9: sentence with two bold runs + its source URL --------------------
$a9Bold1  = "Nitko"
$a9Mid    = " dosad nije "
$a9Bold2  = "imao"
$a9Suffix = " ovako kompleksne uvjete za Schengen«, izjavio je Božinović u intervjuu za Mediaservis"
$ws.Range("A9").Value = $a9Bold1 + $a9Mid + $a9Bold2 + $a9Suffix

$pos = 1
$len = $a9Bold1.Length
$run = $ws.Range("A9").Characters($pos, $len)
$run.Font.Bold = $true
$run.Font.Name = "Calibri"
$run.Font.Size = 12
$pos += $len

$len = $a9Mid.Length
$run = $ws.Range("A9").Characters($pos, $len)
$run.Font.Name = "Calibri"
$run.Font.Size = 12
$pos += $len

$len = $a9Bold2.Length
$run = $ws.Range("A9").Characters($pos, $len)
$run.Font.Bold = $true
$run.Font.Name = "Calibri"
$run.Font.Size = 12
$pos += $len

$len = $a9Suffix.Length
$run = $ws.Range("A9").Characters($pos, $len)
$run.Font.Name = "Calibri"
$run.Font.Size = 12

$ws.Range("B9").Value = "http://novilist.hr/Vijesti/Hrvatska/PITANJE-SCHENGENA-Slovenija-brusi-alate-za-blokadu-Hrvatske?meta_refresh=true"

# --- Row 10: sentence with one bold run, no URL ------------------------------
$a10Prefix = "Župani su u subotu bili gosti Festivala voća "
$a10Bold   = "u mjestu Tavankut"
$a10Suffix = ", a zatim…"
$ws.Range("A10").Value = $a10Prefix + $a10Bold + $a10Suffix

$pos = 1
$len = $a10Prefix.Length
$pos += $len

$len = $a10Bold.Length
$run = $ws.Range("A10").Characters($pos, $len)
$run.Font.Bold = $true
$run.Font.Name = "Calibri"
$run.Font.Size = 12
$pos += $len

$len = $a10Suffix.Length
$run = $ws.Range("A10").Characters($pos, $len)
$run.Font.Name = "Calibri"
$run.Font.Size = 12

# --- Move the active selection to the next empty row, as left by the author -
$ws.Range("A11").Select() | Out-Null
